$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 66693215
$ws.Range("B2").Value = "QUANTIDADE 5 UNIDADES - VOLUME 2 CXS"
